$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.351.31'
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').Value = '3.455.15'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.81'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.69'
$ws.Range('E6').Value = '  +4.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.456.53'
$ws.Range('E8').Value = '  +1.86%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.589'
$ws.Range('E9').Value = '  +10.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.32'
$ws.Range('E10').Value = '  -1.85%  '
$ws.Range('E11').Value = '  +5.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.445'
$ws.Range('E12').Value = '  +2.12%  '
$ws.Range('D13').Value = '4.052.72'
$ws.Range('E13').Value = '  +1.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000195'
$ws.Range('E15').Value = '  +5.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.03'
$ws.Range('E16').Value = '  +7.56%  '
$ws.Range('D17').Value = '64.371.33'
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').Value = '3.407.99'
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.55'
$ws.Range('E20').Value = '  +4.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '388.62'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.22'
$ws.Range('E22').Value = '  -2.37%  '
$ws.Range('E23').Value = '  +2.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.24'
$ws.Range('E24').Value = '  +3.43%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('E26').Value = '  +21.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.58'
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.179'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.20'
$ws.Range('E30').Value = '  +11.35%  '
$ws.Range('E31').Value = '  +9.56%  '
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('E33').Value = '  +1.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.74'
$ws.Range('E34').Value = '  +2.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.12'
$ws.Range('E36').Value = '  +6.08%  '
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.63'
$ws.Range('E38').Value = '  +1.34%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.89'
$ws.Range('E39').Value = '  +0.38%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0776'
$ws.Range('E40').Value = '  +3.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.48'
$ws.Range('E41').Value = '  -0.72%  '
$ws.Range('D42').Value = '2.905.42'
$ws.Range('E42').Value = '  +0.93%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.57'
$ws.Range('E43').Value = '  +6.71%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0320'
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.71'
$ws.Range('E45').Value = '  +3.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.773'
$ws.Range('E46').Value = '  +1.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.88'
$ws.Range('E47').Value = '  +8.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.09'
$ws.Range('E48').Value = '  +2.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.23'
$ws.Range('E49').Value = '  +17.45%  '
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.874'
$ws.Range('E50').Value = '  +7.74%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.108'
$ws.Range('E51').Value = '  +5.09%  '
